$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1. Workbook-level defined names: widen the two named ranges on
#    Formulaire by the 3 newly-inserted columns (AD->AE, AH->AK).
# ------------------------------------------------------------------
$wb.Names.Item(1).RefersTo = "=Formulaire!`$A`$1:`$AE`$1"
$wb.Names.Item(2).RefersTo = "=Formulaire!`$A`$1:`$AK`$1"

# ------------------------------------------------------------------
# 2. Insert "chsta_sprep" as a new column T (it lands just before the
#    existing "chsta_largeurlitmineur" column, shifting everything
#    from T onward one slot to the right).
# ------------------------------------------------------------------
$ws.Columns("T:T").Insert()
$ws.Cells.Item(1, 20).ClearFormats()

# ------------------------------------------------------------------
# 3. Insert "chsta_ombrage" and "chsta_facies" as two new columns,
#    right before "chsta_numphoto" (which, after step 2, lives at
#    column AH).
# ------------------------------------------------------------------
$ws.Columns("AH:AI").Insert()
$ws.Range($ws.Cells.Item(1, 34), $ws.Cells.Item(1, 35)).ClearFormats()

# ------------------------------------------------------------------
# 4. Populate the header text for the 3 new columns. Write them in
#    ombrage / facies / sprep order so the shared-string table gets
#    the same append order as the target workbook.
# ------------------------------------------------------------------
$ws.Cells.Item(1, 34).Value2 = "chsta_ombrage"
$ws.Cells.Item(1, 35).Value2 = "chsta_facies"
$ws.Cells.Item(1, 20).Value2 = "chsta_sprep"

# ------------------------------------------------------------------
# 5. Column widths for the 3 new columns (custom, not best-fit) -
#    reuse widths that already exist elsewhere in the sheet, closest
#    we can land on through the COM ColumnWidth setter.
# ------------------------------------------------------------------
$ws.Columns("T:T").ColumnWidth = 11.75
$ws.Columns("AH:AI").ColumnWidth = 9.6

# ------------------------------------------------------------------
# 6. View state: keep only column A frozen (xSplit = 1) and select
#    T5, matching the saved selection in the target workbook.
# ------------------------------------------------------------------
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("T5").Select()
